$d = $word.ActiveDocument

# 1. Replace the title text
$d.Content.Find.Execute("2.2 - Debate I", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Placeholder - Check Back Later", 2)

# 2. Remove the trailing " " and ":::" runs after the "do not need to be looked at..." sentence
$d.Content.Find.Execute("general edification later. :::", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "general edification later.", 2)
